$wb = $excel.ActiveWorkbook

# ---- Sheet: Weekly ----
$ws1 = $wb.Worksheets.Item("Weekly")
# Extend formatting from the last existing row (93) down into new row 94,
# then set the actual values for the new row.
$ws1.Range("A93:B93").Copy($ws1.Range("A94:B94"))
$ws1.Cells.Item(94,1).Value = 45119
$ws1.Cells.Item(94,2).Value = 5962.440000000001

# ---- Sheet: Resampled2Daily ----
$ws2 = $wb.Worksheets.Item("Resampled2Daily")
# Extend formatting from the last 11 existing rows (638:648) down into the new
# rows 649:659, then set the actual values for each new row.
$ws2.Range("A638:B648").Copy($ws2.Range("A649:B659"))
$ws2.Cells.Item(649,1).Value = 45117
$ws2.Cells.Item(649,2).Value = 6070.89
$ws2.Cells.Item(650,1).Value = 45118
$ws2.Cells.Item(650,2).Value = 6107.075
$ws2.Cells.Item(651,1).Value = 45119
$ws2.Cells.Item(651,2).Value = 5962.440000000001
$ws2.Cells.Item(652,1).Value = 45120
$ws2.Cells.Item(652,2).Value = 6015.154000000001
$ws2.Cells.Item(653,1).Value = 45121
$ws2.Cells.Item(653,2).Value = 6041.809000000001
$ws2.Cells.Item(654,1).Value = 45122
$ws2.Cells.Item(654,2).Value = 6041.809000000001
$ws2.Cells.Item(655,1).Value = 45123
$ws2.Cells.Item(655,2).Value = 6041.809000000001
$ws2.Cells.Item(656,1).Value = 45124
$ws2.Cells.Item(656,2).Value = 6054.264000000001
$ws2.Cells.Item(657,1).Value = 45125
$ws2.Cells.Item(657,2).Value = 6065.724000000001
$ws2.Cells.Item(658,1).Value = 45126
$ws2.Cells.Item(658,2).Value = 6049.782000000001
$ws2.Cells.Item(659,1).Value = 45127
$ws2.Cells.Item(659,2).Value = 6049.782000000001

# ---- Sheet: Daily_TGAData ----
$ws3 = $wb.Worksheets.Item("Daily_TGAData")
# Fix the values for the last three existing rows (corrected data).
$ws3.Cells.Item(646,2).Value = 5971.019
$ws3.Cells.Item(647,2).Value = 5971.019
$ws3.Cells.Item(648,2).Value = 5971.019

# Extend formatting from the last 11 existing rows (638:648) down into the new
# rows 649:659, then set the actual values for each new row.
$ws3.Range("A638:B648").Copy($ws3.Range("A649:B659"))
$ws3.Cells.Item(649,1).Value = 45117
$ws3.Cells.Item(649,2).Value = 5974.554
$ws3.Cells.Item(650,1).Value = 45118
$ws3.Cells.Item(650,2).Value = 5979.732
$ws3.Cells.Item(651,1).Value = 45119
$ws3.Cells.Item(651,2).Value = 5959.359000000001
$ws3.Cells.Item(652,1).Value = 45120
$ws3.Cells.Item(652,2).Value = 6006.791000000001
$ws3.Cells.Item(653,1).Value = 45121
$ws3.Cells.Item(653,2).Value = 6034.524
$ws3.Cells.Item(654,1).Value = 45122
$ws3.Cells.Item(654,2).Value = 6034.524
$ws3.Cells.Item(655,1).Value = 45123
$ws3.Cells.Item(655,2).Value = 6034.524
$ws3.Cells.Item(656,1).Value = 45124
$ws3.Cells.Item(656,2).Value = 6037.184
$ws3.Cells.Item(657,1).Value = 45125
$ws3.Cells.Item(657,2).Value = 6018.542
$ws3.Cells.Item(658,1).Value = 45126
$ws3.Cells.Item(658,2).Value = 6002.6
$ws3.Cells.Item(659,1).Value = 45127
$ws3.Cells.Item(659,2).Value = 6002.6

Write-Host "Edit complete"
